# Update the "times" schedule sheet with the new LA-time (column C) values,
# add a break marker in column D, and move the saved selection to C19.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("times")
$ws.Activate() | Out-Null

# demographic / break row: LA time pushed back 5 minutes
$ws.Range("C10").Value = 1350

# break (row 11): add LA time + a blank-note marker in column D (mirrors D5)
$ws.Range("C11").Value = 1430
$ws.Range("D11").Value = " "

# activity 2 (row 12)
$ws.Range("C12").Value = 1445

# discussion (row 13)
$ws.Range("C13").Value = 1600

# end (row 14)
$ws.Range("C14").Value = 1630

# recap, questions (row 17)
$ws.Range("C17").Value = 800

# intro to projections (row 18)
$ws.Range("C18").Value = 805

# Leave the cursor where the author last left it when saving.
$ws.Range("C19").Select() | Out-Null
